{"js": "// Replace the two-digit multiplication problems' texts with the new values.\n// Each old expression is unique in the document, so a plain text search\n// (matchCase, no wildcards) safely targets exactly one run each.\nconst replacements = [\n  [\"16\u00d785=\", \"28\u00d748=\"],\n  [\"63\u00d798=\", \"85\u00d759=\"],\n  [\"28\u00d799=\", \"49\u00d761=\"],\n  [\"17\u00d779=\", \"72\u00d715=\"],\n  [\"63\u00d736=\", \"46\u00d746=\"],\n  [\"96\u00d766=\", \"29\u00d765=\"],\n  [\"20\u00d787=\", \"15\u00d780=\"],\n  [\"50\u00d755=\", \"69\u00d770=\"],\n  [\"27\u00d764=\", \"28\u00d753=\"],\n  [\"35\u00d722=\", \"51\u00d731=\"],\n  [\"89\u00d763=\", \"66\u00d766=\"],\n  [\"92\u00d778=\", \"66\u00d769=\"],\n  [\"63\u00d737=\", \"34\u00d752=\"],\n  [\"91\u00d767=\", \"18\u00d746=\"],\n  [\"27\u00d731=\", \"42\u00d778=\"],\n  [\"71\u00d730=\", \"71\u00d752=\"],\n  [\"45\u00d795=\", \"18\u00d790=\"],\n  [\"37\u00d728=\", \"22\u00d738=\"],\n  [\"33\u00d789=\", \"68\u00d775=\"],\n  [\"41\u00d772=\", \"64\u00d750=\"],\n  [\"65\u00d711=\", \"73\u00d738=\"],\n  [\"68\u00d757=\", \"66\u00d745=\"],\n  [\"25\u00d720=\", \"14\u00d771=\"],\n  [\"18\u00d778=\", \"18\u00d773=\"],\n  [\"77\u00d768=\", \"93\u00d741=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems' texts with the new values.\n# Each old expression occurs exactly once in the document, so Find/Replace\n# targeting the whole document content safely replaces exactly one run each.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"16\u00d785=\", \"28\u00d748=\"),\n    @(\"63\u00d798=\", \"85\u00d759=\"),\n    @(\"28\u00d799=\", \"49\u00d761=\"),\n    @(\"17\u00d779=\", \"72\u00d715=\"),\n    @(\"63\u00d736=\", \"46\u00d746=\"),\n    @(\"96\u00d766=\", \"29\u00d765=\"),\n    @(\"20\u00d787=\", \"15\u00d780=\"),\n    @(\"50\u00d755=\", \"69\u00d770=\"),\n    @(\"27\u00d764=\", \"28\u00d753=\"),\n    @(\"35\u00d722=\", \"51\u00d731=\"),\n    @(\"89\u00d763=\", \"66\u00d766=\"),\n    @(\"92\u00d778=\", \"66\u00d769=\"),\n    @(\"63\u00d737=\", \"34\u00d752=\"),\n    @(\"91\u00d767=\", \"18\u00d746=\"),\n    @(\"27\u00d731=\", \"42\u00d778=\"),\n    @(\"71\u00d730=\", \"71\u00d752=\"),\n    @(\"45\u00d795=\", \"18\u00d790=\"),\n    @(\"37\u00d728=\", \"22\u00d738=\"),\n    @(\"33\u00d789=\", \"68\u00d775=\"),\n    @(\"41\u00d772=\", \"64\u00d750=\"),\n    @(\"65\u00d711=\", \"73\u00d738=\"),\n    @(\"68\u00d757=\", \"66\u00d745=\"),\n    @(\"25\u00d720=\", \"14\u00d771=\"),\n    @(\"18\u00d778=\", \"18\u00d773=\"),\n    @(\"77\u00d768=\", \"93\u00d741=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    # wdFindContinue = 1, wdReplaceOne = 2\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
